$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the duplicate "Missouri" row that was incorrectly listed under
# the "Sul" (South) region (row 13). Missouri should only appear once,
# under "Centro-Oeste" (Midwest).
$ws.Rows.Item(13).Delete()

# Scroll/selection state to match the author's saved view.
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("B26").Select()
